$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the September 2020 update: new monthly CPI data for rows 15-19
#     (Abr-Ago 2020), matching the existing B:D number-format/style by
#     copying formats from the row above the new data block. ---

$ws.Range("B2:D2").Copy()
$ws.Range("B15:D19").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B15").Value = 310.1243
$ws.Range("C15").Value = 328.7785
$ws.Range("D15").Value = 337.7523

$ws.Range("B16").Value = 314.9087
$ws.Range("C16").Value = 331.0146
$ws.Range("D16").Value = 341.3461

$ws.Range("B17").Value = 321.9738
$ws.Range("C17").Value = 334.4636
$ws.Range("D17").Value = 348.9759

$ws.Range("B18").Value = 328.2014
$ws.Range("C18").Value = 338.7648
$ws.Range("D18").Value = 356.5467

$ws.Range("B19").Value = 337.0632
$ws.Range("C19").Value = 350.5076
$ws.Range("D19").Value = 365.1113

# --- Add new row 20 (Sep-2020 date, same style as the other date cells) ---

$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").Value = 44075

# --- Misc sheet-view bookkeeping the author left behind after finishing ---

$ws.StandardWidth = 11.60546875
$ws.Range("C31").Select()
